$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.100.52"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.879.61"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.117"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.262"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("D13").Value = "1.881.43"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.227"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001102"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06655"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.100"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").Value = "28.140.77"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.267"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").Value = "2.102.09"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "156.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1050"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.651"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.609"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.757"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02455"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06556"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2179"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.224"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6546"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.245"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.914"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6253"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.299"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.685"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.023"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.96%  "
